$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 158-159; this shifts the existing rows 158:190 down to 160:192
# and updates the sheet dimension automatically.
$ws.Rows("158:159").Insert()

# New row 158 ("Especial" quality entry)
$ws.Range("A158").Value = 3
$ws.Range("B158").Value = "Femacal de La Calera"
$ws.Range("C158").Value = "Coquimbo"
$ws.Range("D158").Value = 44543
$ws.Range("E158").Value = 5
$ws.Range("F158").Value = "Fruta"
$ws.Range("G158").Value = 100101
$ws.Range("H158").Value = "Berries"
$ws.Range("I158").Value = 100112025
$ws.Range("J158").Value = "Frutilla"
$ws.Range("K158").Value = "Sin especificar"
$ws.Range("L158").Value = "Especial"
$ws.Range("M158").Value = 68
$ws.Range("N158").Value = 7000
$ws.Range("O158").Value = 7000
$ws.Range("P158").Value = 7000
$ws.Range("Q158").Value = "`$/bandeja 7 kilos"
$ws.Range("R158").Value = "Provincia de Melipilla"
$ws.Range("S158").Value = 1000
$ws.Range("T158").Value = 7

# New row 159 ("Segunda" quality entry)
$ws.Range("A159").Value = 3
$ws.Range("B159").Value = "Femacal de La Calera"
$ws.Range("C159").Value = "Coquimbo"
$ws.Range("D159").Value = 44543
$ws.Range("E159").Value = 5
$ws.Range("F159").Value = "Fruta"
$ws.Range("G159").Value = 100101
$ws.Range("H159").Value = "Berries"
$ws.Range("I159").Value = 100112025
$ws.Range("J159").Value = "Frutilla"
$ws.Range("K159").Value = "Sin especificar"
$ws.Range("L159").Value = "Segunda"
$ws.Range("M159").Value = 50
$ws.Range("N159").Value = 5000
$ws.Range("O159").Value = 5000
$ws.Range("P159").Value = 5000
$ws.Range("Q159").Value = "`$/bandeja 7 kilos"
$ws.Range("R159").Value = "Provincia de Melipilla"
$ws.Range("S159").Value = 714
$ws.Range("T159").Value = 7
